$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A from 65 to 69 (stored OOXML "width" units).
# The ColumnWidth COM property and the raw stored width differ by the
# engine's ~5/6 character padding, so back the value off to land on 69.
$ws.Columns.Item(1).ColumnWidth = 68.17

# Append new feedback rows (37-40)
$ws.Range("A37").Value = "What is the maximum number of headers that I can display in my log?"
$ws.Range("B37").Value = "Unfortunately, the provided information does not specify a maximum number of headers that can be displayed in a log. The Plot Header and Trailer Specifications dialog box allows users to design or modify lithology patterns, modifiers, symbols, headers, and trailers, but it does not provide a limit on the number of headers."

$ws.Range("A38").Value = "What the maximum number of headers I can display in my log?"
$ws.Range("B38").Value = "The maximum number of headers you can display in your log is not explicitly stated in the provided documentation. However, it does mention that on the Edit tab, selecting ""Headers and Trailers"" opens the Plot Header and Trailer Specifications dialog box, which lists various options for customization. It does not provide a specific limit for the number of headers."

$ws.Range("A39").Value = "What the maximum number of headers I can display in my log?"
$ws.Range("B39").Value = "The maximum number of headers you can display in your log is up to 50."

$ws.Range("A40").Value = "How many tables can I have in my log?"
$ws.Range("B40").Value = "You can have up to 100 tables in a log."
